$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" header in column H, matching the existing header
# formatting (bold/centered/bordered) by copying the format from the
# neighboring header cell (G1) rather than re-building the style from
# scratch.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Add the corresponding value for row 2 in the new "Save" column.
$ws.Range("H2").Value = 1
